$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "panel_query_time" style timestamps (column F) on the data sheet ---
$newTimes = @(
    "2021-10-05 14:21:17.956126",
    "2021-10-05 14:21:17.956132",
    "2021-10-05 14:21:17.956135",
    "2021-10-05 14:21:17.956137",
    "2021-10-05 14:21:17.956140",
    "2021-10-05 14:21:17.956142",
    "2021-10-05 14:21:17.956145",
    "2021-10-05 14:21:17.956147",
    "2021-10-05 14:21:17.956149",
    "2021-10-05 14:21:17.956152",
    "2021-10-05 14:21:17.956154",
    "2021-10-05 14:21:17.956156",
    "2021-10-05 14:21:17.956158",
    "2021-10-05 14:21:17.956161",
    "2021-10-05 14:21:17.956163",
    "2021-10-05 14:21:17.956165",
    "2021-10-05 14:21:17.956167",
    "2021-10-05 14:21:17.956169",
    "2021-10-05 14:21:17.956172",
    "2021-10-05 14:21:17.956175"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $newTimes[$i]
}

# --- Add the new "metadata" worksheet right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold/bordered style matching the "data" sheet's header row)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)

$dataSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Intestinal failure"
$metaSheet.Range("C2").Value = 514
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.42"
$metaSheet.Range("D2").ClearFormats()
$metaSheet.Range("E2").Value = "2021-07-12T11:37:33.835303Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:17.953526"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/514/?format=json"

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
